$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.057
$ws.Range("D4").Value = -7.769999999999999
$ws.Range("A7").Value = -20.064
$ws.Range("D12").Value = -7.072
$ws.Range("A16").Value = -22.119
$ws.Range("D18").Value = -8.496
$ws.Range("D19").Value = -8.125
$ws.Range("D20").Value = -7.683999999999999
$ws.Range("A28").Value = -21.933
$ws.Range("A29").Value = -21.325
$ws.Range("D31").Value = -7.850999999999999
$ws.Range("A32").Value = -21.805
$ws.Range("A40").Value = -19.841
$ws.Range("D40").Value = -8.16
$ws.Range("D42").Value = -8.187999999999999
$ws.Range("D47").Value = -7.569
$ws.Range("D48").Value = -7.644
$ws.Range("A52").Value = -21.957
$ws.Range("A57").Value = -22.288
$ws.Range("D63").Value = -7.179
$ws.Range("D64").Value = -7.255
$ws.Range("A66").Value = -21.53
$ws.Range("D76").Value = -7.773000000000001
$ws.Range("D81").Value = -7.540999999999999
$ws.Range("D89").Value = -8.135
$ws.Range("D94").Value = -7.731
$ws.Range("A100").Value = -22.387
